$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column W ("time") with incrementing values starting at 200 for row 3
for ($row = 3; $row -le 255; $row++) {
    $value = 200 + ($row - 3)
    $ws.Cells.Item($row, 23).Value = $value
}

# Update the view: scroll position and selection
$ws.Range("V257").Select()
$excel.ActiveWindow.ScrollRow = 242
$excel.ActiveWindow.ScrollColumn = 2
